$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 & 3: reorder the "Recorded By" email lists (G2, G3) ---
$ws.Range("G2").Value = "Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G3").Value = "asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# --- Class Statistics: Missing Sessions / Pending Sessions counters ---
$ws.Range("L7").Value = 1
$ws.Range("L8").Value = 22

# --- Row 15 (PARASITOLOGY, session 2) flips from Pending -> Not Recorded ---
$ws.Range("I15").Value = "Not Recorded"
$ws.Range("P15").Value = 1
$ws.Range("Q15").Value = 22

# Recolor A15:I15 with the "Not Recorded" (red/pink) status fill - the same
# pink already used for the legend's "Red" swatch (L21) - while keeping the
# row's existing (already-black) font and center alignment untouched.
$a15 = $ws.Range("A15")
$a15.Interior.Color = 12695295
$a15.Interior.PatternColor = 12695295
$a15.Copy()
$ws.Range("B15:I15").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# --- Column I width widened (10 -> 14) ---
$ws.Columns("I").ColumnWidth = 14
